# Apply updated cryptocurrency price/volume figures (cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.261.79"
$ws.Range("E2").Value = "  +1.15%  "
$ws.Range("D3").Value = "1.830.11"
$ws.Range("E3").Value = "  +0.63%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.012"
$ws.Range("E4").Value = "  +0.89%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.04"
$ws.Range("E5").Value = "  +1.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.010"
$ws.Range("E6").Value = "  +0.73%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4728"
$ws.Range("E7").Value = "  +1.52%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3679"
$ws.Range("E8").Value = "  +0.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07433"
$ws.Range("E9").Value = "  +1.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8843"
$ws.Range("E10").Value = "  +1.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.47"
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("D12").Value = "1.913.72"
$ws.Range("E12").Value = "  +5.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07313"
$ws.Range("E13").Value = "  +2.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.417"
$ws.Range("E14").Value = "  +0.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.63"
$ws.Range("E15").Value = "  +2.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.544"
$ws.Range("E16").Value = "  +0.40%  "
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008791"
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("D20").Value = "27.677.57"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.75"
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.283"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.64"
$ws.Range("E23").Value = "  +0.46%  "
$ws.Range("D24").Value = "2.117.57"
$ws.Range("E24").Value = "  +3.25%  "
$ws.Range("E25").Value = "  +0.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.78"
$ws.Range("E26").Value = "  +0.56%  "
$ws.Range("E27").Value = "  +1.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.133"
$ws.Range("E28").Value = "  -0.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.226"
$ws.Range("E29").Value = "  -0.56%  "
$ws.Range("E30").Value = "  -0.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08976"
$ws.Range("E31").Value = "  +0.92%  "
$ws.Range("E32").Value = "  +1.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7457"
$ws.Range("E33").Value = "  -1.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.530"
$ws.Range("E34").Value = "  +0.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.948"
$ws.Range("E35").Value = "  +1.38%  "
$ws.Range("E36").Value = "  +0.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.093"
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05333"
$ws.Range("E38").Value = "  +0.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01951"
$ws.Range("E39").Value = "  +0.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.415"
$ws.Range("E40").Value = "  +2.94%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.951"
$ws.Range("E41").Value = "  -0.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.217"
$ws.Range("E42").Value = "  +0.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5277"
$ws.Range("E43").Value = "  -0.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1653"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.475"
$ws.Range("E45").Value = "  +0.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4899"
$ws.Range("E46").Value = "  +0.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.47"
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "105.00"
$ws.Range("E48").Value = "  +1.50%  "
$ws.Range("E49").Value = "  +0.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.661"
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06303"
$ws.Range("E51").Value = "  +0.04%  "
